# UPDATED TILL - 17-09-2020
# Re-labels the DOB / Applicant Name headers with an input-format hint on a
# second line, widens their columns to fit, and boxes the whole header row
# with a thin border while center-aligning it vertically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text: add the "expected format" hint as a second line ---------
$ws.Range("C1").Value = "DOB" + [char]10 + "YYYY-MM-DD"
$ws.Range("B1").Value = "Applicant Name" + [char]10 + "FIRST_NAME.LAST_NAME"

# --- Column widths: make room for the extra text -----------------------
$ws.Columns.Item(2).ColumnWidth = 32.5
$ws.Columns.Item(3).ColumnWidth = 20.15625

# --- Header row (A1:H1): thin box border + vertically centered ------------
$hdr = $ws.Range("A1:H1")
$hdr.VerticalAlignment = -4108
$hdr.Borders.LineStyle = 1

# --- DOB / Applicant Name headers additionally wrap onto two lines --------
$wrapHdr = $ws.Range("B1:C1")
$wrapHdr.WrapText = $true

# --- Row 1 grows to fit the two-line headers -------------------------------
$ws.Rows.Item(1).RowHeight = 31.5

# --- The stray formatted cell below picks up the same thin border ---------
$ws.Range("F2").Borders.LineStyle = 1

# --- Selection moves to C1 (matches the saved view) ------------------------
$ws.Range("C1").Select()
